$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new product rows (11-13)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "PD Product 10"
$ws.Range("C11").Value = 56
$ws.Range("D11").Value = 100

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "PD Name 11"
$ws.Range("C12").Value = 1212
$ws.Range("D12").Value = 1212

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "asdf"
$ws.Range("C13").Value = 123
$ws.Range("D13").Value = 12

# Update the active selection to match the committed state
$ws.Range("K7").Select()

$wb.Save()
